$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '68.172.56'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '3.496.66'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '595.73'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").Value = '182.71'
$ws.Range("E6").Value = '  +2.95%  '
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  +4.76%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '3.493.18'
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  +7.53%  '
$ws.Range("D11").Value = '7.02'
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("D12").Value = '0.432'
$ws.Range("D13").Value = '4.093.19'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").Value = '32.26'
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '68.147.87'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '0.0000179'
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").Value = '3.491.49'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").Value = '6.23'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '14.21'
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("D21").Value = '396.52'
$ws.Range("E21").Value = '  +2.01%  '
$ws.Range("D22").Value = '8.00'
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("D23").Value = '5.84'
$ws.Range("E23").Value = '  +2.23%  '
$ws.Range("D24").Value = '0.541'
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").Value = '72.31'
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("E27").Value = '  +1.56%  '
$ws.Range("D28").Value = '10.48'
$ws.Range("E28").Value = '  +2.32%  '
$ws.Range("D29").Value = '0.177'
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '6.16'
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").Value = '2.07'
$ws.Range("E33").Value = '  +1.05%  '
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("D35").Value = '7.39'
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").Value = '1.59'
$ws.Range("E37").Value = '  -3.01%  '
$ws.Range("D38").Value = '161.90'
$ws.Range("D39").Value = '0.898'
$ws.Range("E39").Value = '  +3.18%  '
$ws.Range("D40").Value = '2.88'
$ws.Range("E40").Value = '  +11.84%  '
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("D42").Value = '4.72'
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("D43").Value = '6.74'
$ws.Range("E43").Value = '  -3.44%  '
$ws.Range("D44").Value = '26.42'
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("D45").Value = '0.0722'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '26.47'
$ws.Range("E46").Value = '  -3.28%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.758.22'
$ws.Range("E47").Value = '  -2.12%  '
$ws.Range("D48").Value = '41.68'
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("D50").Value = '331.73'
$ws.Range("E50").Value = '  -3.10%  '
$ws.Range("E51").Value = '  -1.94%  '
